# "Generate Report for Handoff"
#
# The localization-status report was regenerated. For the four files that
# are still "Ready for handoff" (rows 4-7 on the zh-cn / de-de / Overview
# sheets), the report run:
#   - bumped their Priority from "low" to "ht"
#   - refreshed the zh-cn "Latest Handoff Datetime" to the new handoff time
#   - refreshed the "Latest HO Xliff Generate Date" / de-de handoff time
#     (these two share the same underlying value)

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Rows 4-7 correspond to the four "Ready for handoff" files on each sheet.
$rows = 4, 5, 6, 7

foreach ($r in $rows) {
    # Priority: low -> ht
    $wsZhCn.Range("E$r").Value = "ht"
    $wsDeDe.Range("E$r").Value = "ht"

    # zh-cn Latest Handoff Datetime
    $wsZhCn.Range("H$r").Value = "2016-08-20 08:41:23"

    # Overview "Latest HO Xliff Generate Date" and de-de Latest Handoff
    # Datetime are the same report timestamp.
    $wsOverview.Range("G$r").Value = "2016-08-20 08:41:27"
    $wsDeDe.Range("H$r").Value = "2016-08-20 08:41:27"
}
